# Apply the 'Updated cryptos list' data refresh (Sun Oct 27 04:50:35 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal-looking string (e.g. '585.47', '1.00')
# must be forced to Text format first, otherwise Excel auto-converts them to
# numbers and drops formatting such as trailing zeros (e.g. '1.00' -> 1).
$forceTextCells = @('D5', 'D15', 'D20', 'D21', 'D24', 'D29', 'D31', 'D35', 'D36', 'D38', 'D42', 'D43', 'D46', 'D48', 'D49', 'D50')
foreach ($cellRef in $forceTextCells) {
    $ws.Range($cellRef).NumberFormat = '@'
}

# Column D / E (and the B/C swap at rows 42-43) updated values
$ws.Range('D2').Value = '67.165.12'
$ws.Range('E2').Value = '  +0.34%  '
$ws.Range('D3').Value = '2.479.91'
$ws.Range('E3').Value = '  +0.84%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '585.47'
$ws.Range('E5').Value = '  +0.91%  '
$ws.Range('E6').Value = '  +3.62%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E8').Value = '  +0.38%  '
$ws.Range('D9').Value = '2.479.06'
$ws.Range('E9').Value = '  +0.81%  '
$ws.Range('E10').Value = '  +3.44%  '
$ws.Range('E11').Value = '  +1.10%  '
$ws.Range('E12').Value = '  +0.87%  '
$ws.Range('E13').Value = '  +0.43%  '
$ws.Range('D14').Value = '2.939.49'
$ws.Range('E14').Value = '  +1.25%  '
$ws.Range('D15').Value = '25.56'
$ws.Range('E15').Value = '  +0.97%  '
$ws.Range('D16').Value = '66.896.64'
$ws.Range('E16').Value = '  +0.58%  '
$ws.Range('E17').Value = '  +1.39%  '
$ws.Range('D18').Value = '2.435.28'
$ws.Range('E18').Value = '  -0.43%  '
$ws.Range('E19').Value = '  +0.03%  '
$ws.Range('D20').Value = '10.99'
$ws.Range('E20').Value = '  -3.04%  '
$ws.Range('D21').Value = '350.60'
$ws.Range('E21').Value = '  -0.84%  '
$ws.Range('E22').Value = '  +0.54%  '
$ws.Range('E23').Value = '  -0.14%  '
$ws.Range('D24').Value = '68.96'
$ws.Range('E24').Value = '  -0.31%  '
$ws.Range('E25').Value = '  +0.54%  '
$ws.Range('E26').Value = '  +3.28%  '
$ws.Range('E27').Value = '  +3.70%  '
$ws.Range('D28').Value = '2.603.57'
$ws.Range('E28').Value = '  +0.86%  '
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.49%  '
$ws.Range('D30').Value = '0.0₃0909'
$ws.Range('E30').Value = '  +1.65%  '
$ws.Range('D31').Value = '509.22'
$ws.Range('E31').Value = '  +0.52%  '
$ws.Range('E32').Value = '  -0.66%  '
$ws.Range('E33').Value = '  +2.49%  '
$ws.Range('E34').Value = '  -0.18%  '
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('D36').Value = '163.37'
$ws.Range('E36').Value = '  +3.40%  '
$ws.Range('E37').Value = '  +2.40%  '
$ws.Range('D38').Value = '18.70'
$ws.Range('E38').Value = '  +0.78%  '
$ws.Range('E39').Value = '  -1.60%  '
$ws.Range('E40').Value = '  -0.43%  '
$ws.Range('E41').Value = '  -0.09%  '
$ws.Range('B42').Value = 'PolygonEcosystemToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D42').Value = '0.330'
$ws.Range('E42').Value = '  +1.29%  '
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').Value = '1.69'
$ws.Range('E43').Value = '  +2.00%  '
$ws.Range('E44').Value = '  +1.39%  '
$ws.Range('E45').Value = '  +3.29%  '
$ws.Range('D46').Value = '143.51'
$ws.Range('E46').Value = '  +1.83%  '
$ws.Range('D47').Value = '0.0₆0262'
$ws.Range('E47').Value = '  +4.06%  '
$ws.Range('D48').Value = '3.49'
$ws.Range('E48').Value = '  +0.62%  '
$ws.Range('D49').Value = '0.516'
$ws.Range('E49').Value = '  +0.49%  '
$ws.Range('D50').Value = '0.0736'
$ws.Range('E50').Value = '  +0.69%  '
$ws.Range('E51').Value = '  -0.32%  '

